# New Test Plans/Try TestCycle
# Update the "Routing Master" sheet with the newly-created Engineering Item
# (Item Number + Id) returned by the latest test run, replacing the
# previous run's values in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-NO6BD"
$ws.Range("D2").Value = "a345f000000uMMoAAM"
